# Source Control Snapshot v2.docx -- re-save refresh.
#
# The upstream commit is a plain re-upload of the document ("Add files via
# upload"): no visible text, formatting, or structural changes. What
# changed under the hood is exactly what Word does whenever a document
# like this is touched/re-saved: the hidden TOC-entry bookmarks that back
# the master document's table of contents (_Toc2473xxxx) get reminted
# with fresh auto-generated names (the _Toc2444xxxx heading bookmarks are
# left alone). Reproduce that bookmark renumbering here.

$d = $word.ActiveDocument

# old hidden bookmark name -> new hidden bookmark name (w:id stays fixed)
$renames = @(
    @{old = "_Toc24731324"; new = "_Toc25148423"},
    @{old = "_Toc24731325"; new = "_Toc25148424"},
    @{old = "_Toc24731326"; new = "_Toc25148425"},
    @{old = "_Toc24731327"; new = "_Toc25148426"},
    @{old = "_Toc24731328"; new = "_Toc25148427"},
    @{old = "_Toc24731329"; new = "_Toc25148428"},
    @{old = "_Toc24731330"; new = "_Toc25148429"}
)

foreach ($pair in $renames) {
    $bm = $d.Bookmarks($pair.old)
    $start = $bm.Start
    $end = $bm.End
    $bm.Delete()
    $rng = $d.Range($start, $end)
    $d.Bookmarks.Add($pair.new, $rng) | Out-Null
}
